# Updated cryptos list on Sat Sep 16 09:42:58 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores numeric-looking values as plain text in the
# source workbook (e.g. "214.96", "26.671.07"). Force text format before
# assigning so Excel doesn't auto-convert them to numbers and drop
# significant trailing zeros / thousands-style dots.
function Set-TextCell($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
}

# --- Row 2: Bitcoin ---
Set-TextCell "D2" "26.671.18"
$ws.Range("E2").Value = "  -0.06%  "

# --- Row 3: Ethereum ---
Set-TextCell "D3" "1.643.17"
$ws.Range("E3").Value = "  +0.68%  "

# --- Row 4: TetherUSD ---
$ws.Range("E4").Value = "  +0.28%  "

# --- Row 5: BNB ---
Set-TextCell "D5" "214.90"
$ws.Range("E5").Value = "  +0.72%  "

# --- Row 6: XRP ---
$ws.Range("E6").Value = "  +0.88%  "

# --- Row 7: USDC ---
$ws.Range("E7").Value = "  +0.28%  "

# --- Row 8: Cardano ---
$ws.Range("E8").Value = "  +0.04%  "

# --- Row 9: Dogecoin ---
Set-TextCell "D9" "0.0628"
$ws.Range("E9").Value = "  +0.87%  "

# --- Row 10: Solana ---
Set-TextCell "D10" "19.27"
$ws.Range("E10").Value = "  +0.11%  "

# --- Row 11: TRON ---
Set-TextCell "D11" "0.0843"
$ws.Range("E11").Value = "  +0.17%  "

# --- Row 12: WrappedliquidstakedEther2.0 ---
Set-TextCell "D12" "1.871.64"
$ws.Range("E12").Value = "  +0.66%  "

# --- Row 13: was WrappedEther, now Polkadot ---
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextCell "D13" "4.21"
$ws.Range("E13").Value = "  +2.84%  "

# --- Row 14: was Polkadot, now WrappedEther ---
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextCell "D14" "1.641.52"
$ws.Range("E14").Value = "  +0.53%  "

# --- Row 15: Polygon ---
$ws.Range("E15").Value = "  +1.16%  "

# --- Row 16: Litecoin ---
$ws.Range("E16").Value = "  +3.17%  "

# --- Row 17: WrappedBTC ---
Set-TextCell "D17" "26.704.77"
$ws.Range("E17").Value = "  +0.15%  "

# --- Row 18: ShibaInu ---
$ws.Range("E18").Value = "  +1.20%  "

# --- Row 19: BitcoinCash ---
Set-TextCell "D19" "216.47"
$ws.Range("E19").Value = "  -1.18%  "

# --- Row 20: Dai ---
$ws.Range("E20").Value = "  +0.22%  "

# --- Row 21: Uniswap ---
$ws.Range("E21").Value = "  +1.47%  "

# --- Row 22: Chainlink ---
$ws.Range("E22").Value = "  +2.31%  "

# --- Row 23: Avalanche ---
$ws.Range("E23").Value = "  +1.94%  "

# --- Row 24: Toncoin ---
Set-TextCell "D24" "2.20"
$ws.Range("E24").Value = "  +12.65%  "

# --- Row 25: Monero ---
Set-TextCell "D25" "145.79"
$ws.Range("E25").Value = "  -1.17%  "

# --- Row 26: BinanceUSD ---
$ws.Range("E26").Value = "  +0.36%  "

# --- Row 27: Stellar ---
$ws.Range("E27").Value = "  -0.54%  "

# --- Row 28: Cosmos ---
$ws.Range("E28").Value = "  +4.64%  "

# --- Row 29: EthereumClassic ---
Set-TextCell "D29" "15.79"
$ws.Range("E29").Value = "  +1.49%  "

# --- Row 30: Hedera ---
Set-TextCell "D30" "0.0517"
$ws.Range("E30").Value = "  +2.34%  "

# --- Row 31: PancakeSwap ---
$ws.Range("E31").Value = "  +0.23%  "

# --- Row 32: Filecoin ---
$ws.Range("E32").Value = "  +2.55%  "

# --- Row 33: InternetComputer(DFINITY) ---
$ws.Range("E33").Value = "  +2.09%  "

# --- Row 34: Maker ---
Set-TextCell "D34" "1.273.13"
$ws.Range("E34").Value = "  +4.72%  "

# --- Row 35: LidoDAOToken ---
$ws.Range("E35").Value = "  +2.22%  "

# --- Row 36: VeChain ---
$ws.Range("E36").Value = "  +6.05%  "

# --- Row 37: HuobiToken ---
$ws.Range("E37").Value = "  +0.19%  "

# --- Row 38: ImmutableX ---
$ws.Range("E38").Value = "  +6.37%  "

# --- Row 39: ARBITRUM ---
Set-TextCell "D39" "0.830"
$ws.Range("E39").Value = "  +2.96%  "

# --- Row 40: PaxDollar ---
$ws.Range("E40").Value = "  +0.31%  "

# --- Row 41: TrustWalletToken ---
Set-TextCell "D41" "0.814"
$ws.Range("E41").Value = "  +2.42%  "

# --- Row 42: MXToken ---
$ws.Range("E42").Value = "  -1.42%  "

# --- Row 43: FraxShare ---
Set-TextCell "D43" "5.46"
$ws.Range("E43").Value = "  +2.14%  "

# --- Row 44: RocketPoolETH ---
Set-TextCell "D44" "1.782.22"
$ws.Range("E44").Value = "  +0.67%  "

# --- Row 45: Quant ---
Set-TextCell "D45" "92.51"

# --- Row 46: Aave ---
Set-TextCell "D46" "59.65"
$ws.Range("E46").Value = "  +8.34%  "

# --- Row 47: RenderToken ---
$ws.Range("E47").Value = "  +2.94%  "

# --- Row 48: was Cronos, now BabyDogeCoin ---
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextCell "D48" "0.0₆0103"
$ws.Range("E48").Value = "  -0.80%  "

# --- Row 49: was EnergySwap, now Cronos ---
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextCell "D49" "0.0516"
$ws.Range("E49").Value = "  +0.79%  "

# --- Row 50: was Algorand, now EnergySwap ---
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell "D50" "7.83"
$ws.Range("E50").Value = "  +3.05%  "

# --- Row 51: was Mantle, now Algorand ---
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextCell "D51" "0.0971"
$ws.Range("E51").Value = "  +3.10%  "
